# Updates the "广州-漫展信息" workbook to the new scrape snapshot.
#
# Summary of the change (per the commit's xml diff):
#   - A new exhibition/event row ("广州·CV谢莹声优见面会专场票·珠三角COMIC WORLD
#     次元世界动漫游戏嘉年华") was scraped and inserted right before the existing
#     "广州·第五人格2.0" row, in both the "展览" sheet (row 29) and the
#     "全部类型" sheet (row 33). All following rows shift down by one.
#   - Several "想去人数" (interest count) values in column F were bumped to
#     reflect newer scrape numbers, across all four sheets.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$cellRef, $value)
    # Force the cell to stay text (Excel otherwise autodetects things like
    # "2024-06-10" as a date and silently converts the stored value).
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

function Apply-FChanges {
    param($ws, $changes)
    foreach ($row in $changes.Keys) {
        $ws.Range("F$row").Value = $changes[$row]
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - A1:I32 -> A1:I33
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# 1) Bump the "want to go" counters for rows that are not affected by the
#    insertion below (rows 1-28 keep their row numbers).
$sheet1Changes = @{
    4  = 115
    6  = 3304
    7  = 981
    8  = 2148
    10 = 1078
    11 = 577
    13 = 1646
    14 = 370
    16 = 29
    18 = 156
    19 = 1524
    20 = 575
    21 = 676
    22 = 563
    23 = 12031
    24 = 12042
    26 = 679
    28 = 9
}
Apply-FChanges $ws1 $sheet1Changes

# 2) Insert the new row at position 29 (shifts old rows 29-32 down to 30-33).
$ws1.Rows.Item(29).Insert()

# Copy formatting (bold/border/centered column A style, etc.) down from the
# row above so the new row matches the rest of the table.
$ws1.Range("A28:I28").Copy()
$ws1.Range("A29:I29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the new row's content.
$ws1.Range("A29").Value = 28
Set-TextValue $ws1 "B29" "2024-06-10"
$ws1.Range("C29").Value = "广州·CV谢莹声优见面会专场票·珠三角COMIC WORLD次元世界动漫游戏嘉年华"
$ws1.Range("D29").Value = "南洲路139号 小洲云文化艺术创意园"
$ws1.Range("E29").Value = "2024.06.10 10:00-06.10 17:00"
$ws1.Range("F29").Value = 0
$ws1.Range("G29").Value = 168
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=85456"
$ws1.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202405/sETdyEqV1715095853125.jpeg"

# 4) Fix up the sequential index column (A) and the "want to go" counters for
#    the rows that shifted down one position.
$ws1.Range("A30").Value = 29
$ws1.Range("F30").Value = 303

$ws1.Range("A31").Value = 30
$ws1.Range("F31").Value = 1886

$ws1.Range("A32").Value = 31
# F32 (重生之道only) keeps its original value of 176 - no change needed.

$ws1.Range("A33").Value = 32
$ws1.Range("F33").Value = 515

# ---------------------------------------------------------------------
# Sheet "演出" (performances) - unaffected by the row insertion, only a few
# "want to go" counters changed.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Changes = @{
    4 = 38
    5 = 117
    6 = 38
}
Apply-FChanges $ws2 $sheet2Changes

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life) - only one counter changed.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Changes = @{
    2 = 71
}
Apply-FChanges $ws3 $sheet3Changes

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) - A1:I40 -> A1:I41, same insertion pattern as
# "展览" but shifted down because this sheet also contains the 演出/本地生活
# rows.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Changes = @{
    3  = 71
    6  = 115
    8  = 3304
    9  = 981
    10 = 2148
    12 = 1078
    13 = 577
    15 = 1646
    16 = 370
    18 = 29
    21 = 38
    22 = 156
    23 = 1524
    24 = 575
    25 = 676
    26 = 563
    27 = 12031
    28 = 12042
    30 = 679
    32 = 9
}
Apply-FChanges $ws4 $sheet4Changes

# Insert the new row at position 33 (shifts old rows 33-40 down to 34-41).
$ws4.Rows.Item(33).Insert()

$ws4.Range("A32:I32").Copy()
$ws4.Range("A33:I33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("A33").Value = 32
Set-TextValue $ws4 "B33" "2024-06-10"
$ws4.Range("C33").Value = "广州·CV谢莹声优见面会专场票·珠三角COMIC WORLD次元世界动漫游戏嘉年华"
$ws4.Range("D33").Value = "南洲路139号 小洲云文化艺术创意园"
$ws4.Range("E33").Value = "2024.06.10 10:00-06.10 17:00"
$ws4.Range("F33").Value = 0
$ws4.Range("G33").Value = 168
$ws4.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=85456"
$ws4.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202405/sETdyEqV1715095853125.jpeg"

$ws4.Range("A34").Value = 33
$ws4.Range("F34").Value = 303

$ws4.Range("A35").Value = 34
$ws4.Range("F35").Value = 1886

$ws4.Range("A36").Value = 35
$ws4.Range("F36").Value = 117

$ws4.Range("A37").Value = 36
$ws4.Range("F37").Value = 38

$ws4.Range("A38").Value = 37
# F38 (重生之道only) keeps its original value of 176 - no change needed.

$ws4.Range("A39").Value = 38
$ws4.Range("F39").Value = 515

$ws4.Range("A40").Value = 39
# F40 (跨越二次元ACG...) keeps its original value of 10 - no change needed.

$ws4.Range("A41").Value = 40
# F41 (孟京辉经典戏剧作品...) keeps its original value of 3 - no change needed.

Write-Host "Applied scrape update (new row inserted + want-to-go counter bumps)."
